# MSO5000 LCR - Clean.xlsx update
#
# - The "I [mA]" column header is relabeled to "I [A]" (current now reported
#   in Amps instead of milliamps), and the I column's numeric values (B2:B20)
#   are rescaled from mA to A (divided by 1000).
# - The active selection on the sheet moves from H1 to H4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the "I" column header from milliamps to amps.
$ws.Range("B1").Value = "I [A]"

# Convert the measured current column from mA to A.
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $mA = $cell.Value()
    $cell.Value = $mA / 1000
}

# Update the active cell selection shown when the workbook is opened.
$ws.Range("H4").Select()
